$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)

# Rename the second sheet (was "Sheet2") to "demoData" and populate it with
# the new demoTestCase data, preserving the cell fill order used originally
# so shared-string indices line up.
$ws2.Name = "demoData"

$ws2.Range("A1").Value = "fName"
$ws2.Range("C1").Value = "email"
$ws2.Range("D1").Value = "company"
$ws2.Range("E1").Value = "phone"

$ws2.Range("A2").Value = "Sunny"
$ws2.Range("B2").Value = "Leone"
$ws2.Range("C2").Value = "sunnyleone@gmail.com"

$ws2.Range("B1").Value = "lName"

$ws2.Range("E2").Value = "Hello"
$ws2.Range("D2").Value = "SunnyLeoneCompany"

# Header row is bold, matching the style already used on sheet1's header.
$ws2.Range("A1:E1").Font.Bold = $true

# Email cell on row 2 gets a mailto hyperlink + the workbook's Hyperlink style.
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:sunnyleone@gmail.com")
$ws2.Range("C2").Style = "Hyperlink"

# Column E (phone) was widened / best-fit to hold its contents.
$ws2.Columns.Item(5).ColumnWidth = 10.14

# Selection on demoData sits on D2, and demoData becomes the active tab.
[void]$ws2.Range("D2").Select()
[void]$ws2.Activate()
